$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ticket #1 (rows 2-17) ---

# "Nº" field
$ws.Range("J4").Value = 3

# "Código:" field
$ws.Range("C7").Value = 3

# "Data:" field -> now plain text instead of a date value
$ws.Range("G7").Value = "14/01/2025"

# "Hora:" field -> now plain text instead of =NOW()
$ws.Range("I7").Value = "10:44"

# "Cliente:" field
$ws.Range("D9").Value = "COCA - COLA"

# "Produto:" field
$ws.Range("D11").Value = "Vidro"

# "Peso Bruto:" field
$ws.Range("D12").Value = 11450

# "Tara:" field
$ws.Range("D13").Value = 8310

# "Peso Liquido:" field
$ws.Range("D14").Value = 3140

# --- Ticket #2 (rows 19-34, mirrors ticket #1) ---

# "Nº" field
$ws.Range("J21").Value = 3

# "Código:" field
$ws.Range("C24").Value = 3

# "Data:" field -> now plain text instead of =TODAY()
$ws.Range("G24").Value = "14/01/2025"

# "Hora:" field -> now plain text instead of =NOW()
$ws.Range("I24").Value = "10:44"

# "Cliente:" field
$ws.Range("D26").Value = "COCA - COLA"

# "Produto:" field
$ws.Range("D28").Value = "Vidro"

# "Peso Bruto:" field
$ws.Range("D29").Value = 11450

# "Tara:" field
$ws.Range("D30").Value = 8310

# "Peso Liquido:" field
$ws.Range("D31").Value = 3140

# Update the active selection to mirror the saved view state
$ws.Activate() | Out-Null
$ws.Range("I24").Select() | Out-Null
